# Update the "Visual" tile-grid sheet to the new pattern.
$wb = $excel.ActiveWorkbook
$visual = $wb.Worksheets.Item("Visual")

$rows = @(
    "11111111",
    "10000001",
    "10100101",
    "10000001",
    "10100101",
    "10011001",
    "10000001",
    "11111111"
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $line = $rows[$r]
    for ($c = 0; $c -lt $line.Length; $c++) {
        $cell = $visual.Cells.Item($r + 1, $c + 1)
        if ($line[$c] -eq '1') {
            $cell.Value = 1
        } else {
            $cell.ClearContents()
        }
    }
}

# Move the active selection on the Visual sheet to match the saved view.
$visual.Activate()
$visual.Range("F6").Select()

$wb.Application.Calculate()
